$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Jordan Network table (rows 2-4): add validation/test result for K3 ---
$ws.Range("K3").Value = " 96.45/93.65"

# --- Elman Keras table (rows 13-18): add baseline row (15) numbering and
#     results for run no. 3 (rows 16-18 / A17 "3") ---
$ws.Range("A15").Value = "1"
$ws.Range("A16").Value = "2"
$ws.Range("A17").Value = "3"
$ws.Range("H17").Value = "95.64/92.22"

$ws.Range("C18").Value = "0.1"
$ws.Range("D18").Value = "200"
$ws.Range("F18").Value = "200"
$ws.Range("G18").Value = "250"
